$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order #5 (row 56) and Order #6 (row 57) have been manually finalized
# ("Auto Done"). Their Min/Max/Middle Wavelength columns (D/E/F) become
# fixed manual values instead of the shared formulas, the Range column
# (G) becomes a simple (non-shared) difference formula, and column K is
# marked "Auto Done".

# Row 56 - Order #5
$ws.Range("D56").Value = 3918.2
$ws.Range("E56").Value = 3983.7
$ws.Range("F56").Value = 3951.6
$ws.Range("G56").Formula = "=E56-D56"
$ws.Range("K56").Value = "Auto Done"

# Row 57 - Order #6
$ws.Range("D57").Value = 3963.9
$ws.Range("E57").Value = 4030.1
$ws.Range("F57").Value = 3997.6
$ws.Range("G57").Formula = "=E57-D57"
$ws.Range("K57").Value = "Auto Done"

# Leave the selection on the last-edited cell, matching the author's session.
$ws.Range("K57").Select()
